# Insert a new weekly data row (row 94) into the "Zanahoria" price table.
# This pushes the existing rows 94..204 down to 95..205, so the table
# grows from A1:R204 to A1:R205.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 94 (shifts rows 94..204 down to 95..205).
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly record.
$ws.Range("A94").Value = 5
$ws.Range("B94").Value = "Macroferia Regional de Talca"
$ws.Range("C94").Value = "Maule"
$ws.Range("D94").Value = 44494
$ws.Range("E94").Value = 7
$ws.Range("F94").Value = 100114013
$ws.Range("G94").Value = "Zanahoria"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 300
$ws.Range("K94").Value = 7000
$ws.Range("L94").Value = 7000
$ws.Range("M94").Value = 7000
$ws.Range("N94").Value = "$/saco 20 kilos"
$ws.Range("O94").Value = "Región de Ñuble"
$ws.Range("P94").Value = 350
$ws.Range("Q94").Value = 20
$ws.Range("R94").Value = "Hortaliza"
